# The workbook's "Data" sheet has a row (row 6) whose "post" (h2o/EtOH)
# measurements were never recorded -- only the "pre" values exist for that
# row. The sheet's xlsx parser needs to account for rows that have "pre"
# data but no "post" data, so the stray zero/placeholder values that were
# sitting in F6 (post h2o) and G6 (post EtOH) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the "post" values for row 6 - this row only has "pre" data.
$ws.Range("F6:G6").ClearContents()

# Move/update the active selection to reflect where the edit was made.
$ws.Range("F6").Select()
